$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Parameters": rename parameter labels to camelCase, clear the
# now-unused description column (C), and move the selection to C8.
# (Edited first so new shared strings land in the same order the original
# author created them in.)
# ---------------------------------------------------------------------------
$wsParams = $wb.Worksheets.Item("Parameters")
$wsParams.Range("A2").Value = "firstRowOfData"
$wsParams.Range("A3").Value = "mockColumn"
$wsParams.Range("A4").Value = "deltaIndicatorColumn"
$wsParams.Range("A5").Value = "concatenatedKeyColumn"
$wsParams.Range("A6").Value = "startColumnCheckData"
$wsParams.Range("A7").Value = "headerRow"
$wsParams.Range("A10").Value = "EnableMockNumberCheck"

$wsParams.Range("C2:C10").ClearContents()

# Row 9 had a custom (taller) row height to fit the long description text that
# used to live in C9; now that the cell is empty, let it size back down.
$wsParams.Rows.Item(9).AutoFit()

[void]$wsParams.Range("C8").Select()

# ---------------------------------------------------------------------------
# Sheet "File List": update Source Path / Source File Name, and make this the
# active / selected sheet with B5 selected.
# ---------------------------------------------------------------------------
$wsFileList = $wb.Worksheets.Item("File List")
$wsFileList.Range("A2").Value = "C:\Users\j.a.vorathammaporn\OneDrive - Accenture\Desktop\PTT-WorkSpace\SandBox\DeltaTest\Python_Result\"
$wsFileList.Range("B2").Value = "BeforeDelta_Python.xlsx"

# ---------------------------------------------------------------------------
# Make "File List" the active sheet/tab, with B5 selected (this also clears
# "Run Program"'s tabSelected flag).
# ---------------------------------------------------------------------------
[void]$wsFileList.Activate()
[void]$wsFileList.Range("B5").Select()
